$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the summary text in A1 with the new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.77 = 51628.8 pesos`n✅ 51628.8 pesos = 12.78 = 977.51 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate values (N10, O10, O12) ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 78.309
$ws2.Range("O10").Value = 4043
$ws2.Range("O12").Value = 76.491
